# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1) Update the "Date" metadata value on the Metadata sheet.
# 2) Swap the two mapping columns (AK = "Mapping: RIM Mapping" and
#    AL = "Mapping: Specification metier vers l'extension ROR
#    HealthcareServiceContact") on the Elements sheet, including the
#    header row, all the data rows and the column widths, so that the
#    business-mapping column now comes before the RIM-mapping column.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 (Date) -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2024-03-19T13:17:15+00:00"

# --- 2) Elements sheet: swap columns AK (37) and AL (38) -------------------
$ws = $wb.Worksheets.Item("Elements")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $ws.Cells.Item($r, 37)
    $alCell = $ws.Cells.Item($r, 38)

    $akValue = $akCell.Value2
    $alValue = $alCell.Value2

    $akCell.Value = $alValue
    $alCell.Value = $akValue
}

# Swap the column widths too (AK was narrow/24.98, AL was wide/81.95;
# after the content swap AK becomes the wide column and AL the narrow one).
$ws.Columns.Item(37).ColumnWidth = 81.09
$ws.Columns.Item(38).ColumnWidth = 24.09
